$d = $word.ActiveDocument
$r = $d.Content
$r.Collapse(0)
$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:spacing w:before="0" w:after="160"/><w:ind w:right="0" w:hanging="0"/><w:contextualSpacing/><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:eastAsia="Calibri" w:cs="Tahoma"/><w:b w:val="false"/><w:b w:val="false"/><w:bCs w:val="false"/><w:color w:val="auto"/><w:kern w:val="0"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="en-IE" w:eastAsia="en-US" w:bidi="ar-SA"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:eastAsia="Calibri" w:cs="Tahoma"/><w:b w:val="false"/><w:bCs w:val="false"/><w:color w:val="auto"/><w:kern w:val="0"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="en-IE" w:eastAsia="en-US" w:bidi="ar-SA"/></w:rPr></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:spacing w:before="0" w:after="160"/><w:ind w:right="0" w:hanging="0"/><w:contextualSpacing/><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:eastAsia="Calibri" w:cs="Tahoma"/><w:b w:val="false"/><w:b w:val="false"/><w:bCs w:val="false"/><w:color w:val="auto"/><w:kern w:val="0"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="en-IE" w:eastAsia="en-US" w:bidi="ar-SA"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:eastAsia="Calibri" w:cs="Tahoma"/><w:b w:val="false"/><w:bCs w:val="false"/><w:color w:val="auto"/><w:kern w:val="0"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="en-IE" w:eastAsia="en-US" w:bidi="ar-SA"/></w:rPr></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Normal"/><w:spacing w:before="0" w:after="160"/><w:ind w:left="0" w:right="0" w:hanging="0"/><w:contextualSpacing/><w:rPr/></w:pPr><w:r><w:rPr/><w:t>Week of</w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsia="Calibri" w:cs="Tahoma"/><w:color w:val="auto"/><w:kern w:val="0"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="en-IE" w:eastAsia="en-US" w:bidi="ar-SA"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsia="Calibri" w:cs="Tahoma"/><w:color w:val="auto"/><w:kern w:val="0"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="en-IE" w:eastAsia="en-US" w:bidi="ar-SA"/></w:rPr><w:t>09</w:t></w:r><w:r><w:rPr/><w:t>/0</w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsia="Calibri" w:cs="Tahoma"/><w:color w:val="auto"/><w:kern w:val="0"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="en-IE" w:eastAsia="en-US" w:bidi="ar-SA"/></w:rPr><w:t>5</w:t></w:r><w:r><w:rPr/><w:t>/21</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="6"/></w:numPr><w:spacing w:before="0" w:after="160"/><w:ind w:left="720" w:right="0" w:hanging="0"/><w:contextualSpacing/><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:eastAsia="Calibri" w:cs="Tahoma"/><w:b w:val="false"/><w:b w:val="false"/><w:bCs w:val="false"/><w:color w:val="auto"/><w:kern w:val="0"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="en-IE" w:eastAsia="en-US" w:bidi="ar-SA"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:eastAsia="Calibri" w:cs="Tahoma"/><w:b w:val="false"/><w:bCs w:val="false"/><w:color w:val="auto"/><w:kern w:val="0"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="en-IE" w:eastAsia="en-US" w:bidi="ar-SA"/></w:rPr><w:t xml:space="preserve">Project  video + poster </w:t></w:r></w:p>'
$r.InsertXML($xml)
